# The scraper added two new columns of player bio data ("height" and
# "weight") right before the existing "fantasy points" column, which
# shifts two columns to the right (E -> G) to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at E:F. This pushes the current column E
# ("fantasy points", header + 16 data rows) over to column G, preserving
# its header text/style and all of its values.
$ws.Range("E:F").Insert()

# Populate the newly inserted column E with the "height" data.
$ws.Range("E1").Value = "height"
$ws.Range("E2:E17").Value = 6.333333333333333

# Populate the newly inserted column F with the "weight" data.
$ws.Range("F1").Value = "weight"
$ws.Range("F2:F17").Value = 246
